# Auto-generated Excel COM-interop script
# Applies numeric data updates to the per-class profit tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as captured by the scheduled market-data refresh runner.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 24000
$ws.Range("J7").Value = 24000
$ws.Range("L7").Value = 24000
$ws.Range("N7").Value = -24224

$ws.Range("H10").Value = 17500
$ws.Range("J10").Value = 17500
$ws.Range("L10").Value = 17500
$ws.Range("N10").Value = -18086

$ws.Range("H13").Value = 20833.666
$ws.Range("J13").Value = 20833.666
$ws.Range("L13").Value = 20833.666
$ws.Range("N13").Value = -21171.666

$ws.Range("H14").Value = 24000
$ws.Range("J14").Value = 24000
$ws.Range("L14").Value = 24000
$ws.Range("N14").Value = -24382

$ws.Range("H16").Value = 7300
$ws.Range("I16").Value = 1500
$ws.Range("J16").Value = 8750
$ws.Range("K16").Value = 1500
$ws.Range("L16").Value = 8750
$ws.Range("M16").Value = -1270
$ws.Range("N16").Value = -9210

$ws.Range("H54").Value = 19297.857
$ws.Range("J54").Value = 19297.857
$ws.Range("L54").Value = 19297.857
$ws.Range("N54").Value = -20269.857

$ws.Range("H58").Value = 23988.37
$ws.Range("J58").Value = 26239.643
$ws.Range("L58").Value = 78718.929
$ws.Range("N58").Value = -79018.929

$ws.Range("H69").Value = 3519
$ws.Range("I69").Value = 2995.3333
$ws.Range("J69").Value = 4042.6667
$ws.Range("K69").Value = 8985.999899999999
$ws.Range("L69").Value = 12128.0001
$ws.Range("M69").Value = -8111.999899999999
$ws.Range("N69").Value = -13876.0001

$ws.Range("H72").Value = 3519
$ws.Range("I72").Value = 2995.3333
$ws.Range("J72").Value = 4042.6667
$ws.Range("K72").Value = 26957.9997
$ws.Range("L72").Value = 36384.0003
$ws.Range("M72").Value = -22589.9997
$ws.Range("N72").Value = -45120.0003

$ws.Range("H96").Value = 847.8125
$ws.Range("I96").Value = 770.7
$ws.Range("K96").Value = 2312.1
$ws.Range("M96").Value = -939.1000000000004

$ws.Range("H135").Value = 1171.9395
$ws.Range("I135").Value = 845.5357
$ws.Range("J135").Value = 2999.8
$ws.Range("K135").Value = 7609.821300000001
$ws.Range("L135").Value = 26998.2
$ws.Range("M135").Value = -5074.821300000001
$ws.Range("N135").Value = -32068.2

$ws.Range("H137").Value = 3321.5
$ws.Range("I137").Value = 3313.2693
$ws.Range("K137").Value = 9939.8079
$ws.Range("M137").Value = -7389.8079

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1372.8529
$ws.Range("I45").Value = 1057.4231
$ws.Range("J45").Value = 2398
$ws.Range("K45").Value = 1057.4231
$ws.Range("L45").Value = 2398
$ws.Range("M45").Value = -680.4231
$ws.Range("N45").Value = -3152

$ws.Range("H51").Value = 45734.43
$ws.Range("J51").Value = 45734.43
$ws.Range("L51").Value = 45734.43
$ws.Range("N51").Value = -47246.43

$ws.Range("H74").Value = 874
$ws.Range("I74").Value = 812.24
$ws.Range("K74").Value = 812.24
$ws.Range("M74").Value = 61.75999999999999

$ws.Range("H77").Value = 874
$ws.Range("I77").Value = 812.24
$ws.Range("K77").Value = 4061.2
$ws.Range("M77").Value = 306.8000000000002

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H141").Value = 39800
$ws.Range("J141").Value = 39800
$ws.Range("L141").Value = 39800
$ws.Range("N141").Value = -50160

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 258.33334
$ws.Range("I12").Value = 258.33334
$ws.Range("K12").Value = 258.33334
$ws.Range("M12").Value = -90.33334000000002

$ws.Range("H22").Value = 185.55556
$ws.Range("I22").Value = 95.71429000000001
$ws.Range("K22").Value = 95.71429000000001
$ws.Range("M22").Value = 77.28570999999999

$ws.Range("H35").Value = 20074
$ws.Range("J35").Value = 20074
$ws.Range("L35").Value = 20074
$ws.Range("N35").Value = -20694

$ws.Range("H133").Value = 14390
$ws.Range("J133").Value = 14390
$ws.Range("L133").Value = 14390
$ws.Range("N133").Value = -24510

$ws.Range("H134").Value = 1904.925
$ws.Range("I134").Value = 1399.9117
$ws.Range("J134").Value = 4766.6665
$ws.Range("K134").Value = 4199.7351
$ws.Range("L134").Value = 14299.9995
$ws.Range("M134").Value = -1664.7351
$ws.Range("N134").Value = -19369.9995

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 52551.5
$ws.Range("J11").Value = 69668.664
$ws.Range("L11").Value = 69668.664
$ws.Range("N11").Value = -69948.664

$ws.Range("H21").Value = 19295.416
$ws.Range("I21").Value = 1000
$ws.Range("J21").Value = 44909
$ws.Range("K21").Value = 1000
$ws.Range("L21").Value = 44909
$ws.Range("M21").Value = -765
$ws.Range("N21").Value = -45379

$ws.Range("H31").Value = 4138.263
$ws.Range("I31").Value = 3166.8572
$ws.Range("J31").Value = 5338.2354
$ws.Range("K31").Value = 3166.8572
$ws.Range("L31").Value = 5338.2354
$ws.Range("M31").Value = -2871.8572
$ws.Range("N31").Value = -5928.2354

$ws.Range("H34").Value = 4138.263
$ws.Range("I34").Value = 3166.8572
$ws.Range("J34").Value = 5338.2354
$ws.Range("K34").Value = 3166.8572
$ws.Range("L34").Value = 5338.2354
$ws.Range("M34").Value = -2964.8572
$ws.Range("N34").Value = -5742.2354

$ws.Range("H52").Value = 29933.334
$ws.Range("J52").Value = 29933.334
$ws.Range("L52").Value = 29933.334
$ws.Range("N52").Value = -30521.334

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 2960
$ws.Range("J19").Value = 2950
$ws.Range("L19").Value = 8850
$ws.Range("N19").Value = -9198

$ws.Range("H44").Value = 1021.75
$ws.Range("I44").Value = 482.2
$ws.Range("J44").Value = 1407.1428
$ws.Range("K44").Value = 1446.6
$ws.Range("L44").Value = 4221.428400000001
$ws.Range("M44").Value = -1048.6
$ws.Range("N44").Value = -5017.428400000001

$ws.Range("H64").Value = 2375.75
$ws.Range("J64").Value = 2999.8
$ws.Range("L64").Value = 8999.400000000001
$ws.Range("N64").Value = -9539.400000000001

$ws.Range("H67").Value = 2375.75
$ws.Range("J67").Value = 2999.8
$ws.Range("L67").Value = 8999.400000000001
$ws.Range("N67").Value = -10871.4

$ws.Range("H131").Value = 1317.4
$ws.Range("I131").Value = 2081.4285
$ws.Range("J131").Value = 1020.2778
$ws.Range("K131").Value = 6244.2855
$ws.Range("L131").Value = 3060.8334
$ws.Range("M131").Value = -1204.2855
$ws.Range("N131").Value = -13140.8334

$ws.Range("H136").Value = 2357.8696
$ws.Range("I136").Value = 1509.9333
$ws.Range("J136").Value = 3947.75
$ws.Range("K136").Value = 4529.7999
$ws.Range("L136").Value = 11843.25
$ws.Range("M136").Value = 570.2001
$ws.Range("N136").Value = -22043.25

$ws.Range("H140").Value = 7940122.5
$ws.Range("I140").Value = 33334144
$ws.Range("J140").Value = 4490.625
$ws.Range("K140").Value = 100002432
$ws.Range("L140").Value = 13471.875
$ws.Range("M140").Value = -99997252
$ws.Range("N140").Value = -23831.875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2801.9583
$ws.Range("I126").Value = 1583.2
$ws.Range("J126").Value = 3672.5
$ws.Range("K126").Value = 4749.6
$ws.Range("L126").Value = 11017.5
$ws.Range("M126").Value = -2279.6
$ws.Range("N126").Value = -15957.5

$ws.Range("H132").Value = 2894.551
$ws.Range("I132").Value = 2387.0833
$ws.Range("J132").Value = 4299.846
$ws.Range("K132").Value = 7161.249899999999
$ws.Range("L132").Value = 12899.538
$ws.Range("M132").Value = -4631.249899999999
$ws.Range("N132").Value = -17959.538

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H51").Value = 26750.666
$ws.Range("J51").Value = 26750.666
$ws.Range("L51").Value = 26750.666
$ws.Range("N51").Value = -27706.666

$ws.Range("H61").Value = 58826996
$ws.Range("I61").Value = 111113016
$ws.Range("J61").Value = 5225
$ws.Range("K61").Value = 111113016
$ws.Range("L61").Value = 5225
$ws.Range("M61").Value = -111112814
$ws.Range("N61").Value = -5629

$ws.Range("H113").Value = 58826996
$ws.Range("I113").Value = 111113016
$ws.Range("J113").Value = 5225
$ws.Range("K113").Value = 111113016
$ws.Range("L113").Value = 5225
$ws.Range("M113").Value = -111110846
$ws.Range("N113").Value = -9565

$ws.Range("H132").Value = 3315.9688
$ws.Range("I132").Value = 2042.3572
$ws.Range("J132").Value = 4306.5557
$ws.Range("K132").Value = 6127.071599999999
$ws.Range("L132").Value = 12919.6671
$ws.Range("M132").Value = -3597.071599999999
$ws.Range("N132").Value = -17979.6671

$ws.Range("H133").Value = 29700
$ws.Range("J133").Value = 29700
$ws.Range("L133").Value = 29700
$ws.Range("N133").Value = -34760

$ws.Range("H135").Value = 29678.625
$ws.Range("J135").Value = 29678.625
$ws.Range("L135").Value = 29678.625
$ws.Range("N135").Value = -39818.625

$ws.Range("H137").Value = 29587.777
$ws.Range("J137").Value = 29587.777
$ws.Range("L137").Value = 29587.777
$ws.Range("N137").Value = -39787.777

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 40000
$ws.Range("J138").Value = 40000
$ws.Range("L138").Value = 40000
$ws.Range("N138").Value = -50280

$ws.Range("H139").Value = 40000
$ws.Range("J139").Value = 40000
$ws.Range("L139").Value = 40000
$ws.Range("N139").Value = -50280
